$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-text edits (preserve rich-text runs via Characters) ---

# A8: "Volume 29   Number  44" -> "...45"  (chars 21-22 are "44")
$ws.Range("A8").Characters(21, 2).Text = "45"

# C9: "Report Covering the Week  10/31/2022  Through  11/6/2022"
#     -> "...11/7/2022  Through  11/13/2022"
$ws.Range("C9").Characters(27, 10).Text = "11/7/2022"
$ws.Range("C9").Characters(47, 9).Text = "11/13/2022"

# --- Row 16: D16/E16 switch from text placeholders ("0"/"***.*") to numeric cells ---
# Copy number-format/style from sibling numeric cells, then set the value.
$ws.Range("C16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 3
$ws.Range("K16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = 133.333333333333

# --- Remaining numeric cell updates ---

# Row 15
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = 5.555555555555
$ws.Range("L15").Value = 46.153846153846
$ws.Range("M15").Value = 46.153846153846
$ws.Range("N15").Value = -17.391304347826

# Row 16
$ws.Range("C16").Value = 7
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 30.76923076923
$ws.Range("I16").Value = 139
$ws.Range("J16").Value = 116
$ws.Range("K16").Value = 19.827586206896
$ws.Range("L16").Value = 19.827586206896
$ws.Range("M16").Value = 15.833333333333
$ws.Range("N16").Value = -88.464730290456

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("I17").Value = 161
$ws.Range("J17").Value = 129
$ws.Range("K17").Value = 24.806201550387
$ws.Range("L17").Value = 46.363636363636
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = -63.325740318906

# Row 18
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -6.666666666666
$ws.Range("I18").Value = 194
$ws.Range("J18").Value = 177
$ws.Range("K18").Value = 9.604519774011
$ws.Range("L18").Value = -14.53744493392
$ws.Range("M18").Value = 13.45029239766
$ws.Range("N18").Value = -89.553042541734

# Row 19
$ws.Range("C19").Value = 37
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = 12.121212121212
$ws.Range("F19").Value = 179
$ws.Range("G19").Value = 133
$ws.Range("H19").Value = 34.586466165413
$ws.Range("I19").Value = 1713
$ws.Range("J19").Value = 1033
$ws.Range("K19").Value = 65.827686350435
$ws.Range("L19").Value = 80.887011615628
$ws.Range("M19").Value = 15.353535353535
$ws.Range("N19").Value = -73.658311548516

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 7
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 107
$ws.Range("J20").Value = 57
$ws.Range("K20").Value = 87.719298245614
$ws.Range("L20").Value = 81.355932203389
$ws.Range("M20").Value = 160.975609756098
$ws.Range("N20").Value = -75.345622119815

# Row 21
$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = 15.555555555555
$ws.Range("F21").Value = 229
$ws.Range("G21").Value = 174
$ws.Range("H21").Value = 31.609195402298
$ws.Range("I21").Value = 2334
$ws.Range("J21").Value = 1533
$ws.Range("K21").Value = 52.25048923679
$ws.Range("L21").Value = 58.130081300813
$ws.Range("M21").Value = 19.815195071868
$ws.Range("N21").Value = -77.711993888464

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -33.333333333333
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 20
$ws.Range("I22").Value = 59
$ws.Range("J22").Value = 45
$ws.Range("K22").Value = 31.111111111111
$ws.Range("L22").Value = -10.60606060606
$ws.Range("M22").Value = 13.461538461538

# Row 23
$ws.Range("M23").Value = -60

# Row 24
$ws.Range("C24").Value = 72
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 282
$ws.Range("G24").Value = 217
$ws.Range("H24").Value = 29.953917050691
$ws.Range("I24").Value = 2381
$ws.Range("J24").Value = 1700
$ws.Range("K24").Value = 40.058823529411
$ws.Range("L24").Value = 91.70692431562
$ws.Range("M24").Value = 42.064439140811

# Row 25
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = 17.948717948717
$ws.Range("I25").Value = 503
$ws.Range("J25").Value = 388
$ws.Range("K25").Value = 29.639175257732
$ws.Range("L25").Value = 82.90909090909
$ws.Range("M25").Value = 27.020202020202

# Row 26
$ws.Range("I26").Value = 32
$ws.Range("K26").Value = 6.666666666666
$ws.Range("L26").Value = 68.421052631578

# Row 27
$ws.Range("C27").Value = 3
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 87
$ws.Range("K27").Value = 17.567567567567
$ws.Range("L27").Value = 58.181818181818

# Row 30
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 50

Write-Output "edits applied"